$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.472.00'
$ws.Range("D3").Value = '1.669.46'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '237.69'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.4796'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").Value = '0.2633'
$ws.Range("E8").Value = '  +1.12%  '
$ws.Range("D9").Value = '0.06179'
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("D10").Value = '0.07031'
$ws.Range("E10").Value = '  -2.24%  '
$ws.Range("D11").Value = '1.668.58'
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").Value = '14.85'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '0.5894'
$ws.Range("E13").Value = '  -5.03%  '
$ws.Range("D14").Value = '4.371'
$ws.Range("E14").Value = '  -3.36%  '
$ws.Range("D15").Value = '75.00'
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '0.9999'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '25.468.81'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("D19").Value = '0.000006757'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = '11.46'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '1.881.17'
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = '4.456'
$ws.Range("E22").Value = '  -1.46%  '
$ws.Range("D23").Value = '8.738'
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("D24").Value = '5.285'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '136.92'
$ws.Range("E25").Value = '  +3.74%  '
$ws.Range("D26").Value = '15.05'
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").Value = '1.389'
$ws.Range("E27").Value = '  -0.65%  '
$ws.Range("D28").Value = '1.723'
$ws.Range("E28").Value = '  +3.05%  '
$ws.Range("D29").Value = '105.09'
$ws.Range("E29").Value = '  +2.13%  '
$ws.Range("D30").Value = '3.949'
$ws.Range("E30").Value = '  +5.14%  '
$ws.Range("D31").Value = '0.07809'
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("D32").Value = '3.647'
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '0.04224'
$ws.Range("E34").Value = '  -6.24%  '
$ws.Range("D35").Value = '2.604'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").Value = '0.6094'
$ws.Range("E36").Value = '  +4.77%  '
$ws.Range("D37").Value = '0.9486'
$ws.Range("E37").Value = '  +1.63%  '
$ws.Range("D38").Value = '2.599'
$ws.Range("E38").Value = '  +0.90%  '
$ws.Range("D39").Value = '0.8564'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("D40").Value = '0.9998'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").Value = '1.852'
$ws.Range("E41").Value = '  +1.92%  '
$ws.Range("D42").Value = '0.01477'
$ws.Range("E42").Value = '  -5.69%  '
$ws.Range("D43").Value = '96.25'
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("D44").Value = '0.3770'
$ws.Range("E44").Value = '  +1.39%  '
$ws.Range("D45").Value = '4.840'
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").Value = '0.1120'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").Value = '6.209'
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").Value = '0.05252'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  +0.15%  '
